# Auto-generated from the commit diff: refresh cached market-board values
# in the Leve-profit worksheets (data-only cells, no formulas in this workbook).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 28575548
$ws.Range("I106").Value = 75005816
$ws.Range("J106").Value = 3076.923
$ws.Range("K106").Value = 75005816
$ws.Range("L106").Value = 3076.923
$ws.Range("M106").Value = -75005185
$ws.Range("N106").Value = -4338.923
$ws.Range("H129").Value = 1551.9584
$ws.Range("I129").Value = 375
$ws.Range("J129").Value = 1658.9546
$ws.Range("K129").Value = 1125
$ws.Range("L129").Value = 4976.8638
$ws.Range("M129").Value = 3875
$ws.Range("N129").Value = -14976.8638
$ws.Range("H135").Value = 17858716
$ws.Range("I135").Value = 23810440
$ws.Range("K135").Value = 214293960
$ws.Range("M135").Value = -214291425
$ws.Range("H137").Value = 3335587.5
$ws.Range("I137").Value = 5210106.5
$ws.Range("J137").Value = 3108.889
$ws.Range("K137").Value = 15630319.5
$ws.Range("L137").Value = 9326.667000000001
$ws.Range("M137").Value = -15627769.5
$ws.Range("N137").Value = -14426.667
$ws.Range("H138").Value = 4169.642
$ws.Range("I138").Value = 3640.95
$ws.Range("J138").Value = 4342.9834
$ws.Range("K138").Value = 10922.85
$ws.Range("L138").Value = 13028.9502
$ws.Range("M138").Value = -5782.849999999999
$ws.Range("N138").Value = -23308.9502

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11114640
$ws.Range("I61").Value = 18520834
$ws.Range("K61").Value = 18520834
$ws.Range("M61").Value = -18520622
$ws.Range("H136").Value = 11114640
$ws.Range("I136").Value = 18520834
$ws.Range("K136").Value = 55562502
$ws.Range("M136").Value = -55559952

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()  # was -80656
$ws.Range("H86").Value = 1938.8462
$ws.Range("I86").Value = 1876.4
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 1876.4
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -753.4000000000001
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 1938.8462
$ws.Range("I89").Value = 1876.4
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 9382
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -3766
$ws.Range("N89").Value = -28732

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2198.2104
$ws.Range("I16").Value = 2349.0833
$ws.Range("K16").Value = 2349.0833
$ws.Range("M16").Value = -2062.0833
$ws.Range("H31").Value = 8329.884
$ws.Range("J31").Value = 8329.884
$ws.Range("L31").Value = 8329.884
$ws.Range("N31").Value = -8919.884
$ws.Range("H34").Value = 8329.884
$ws.Range("J34").Value = 8329.884
$ws.Range("L34").Value = 8329.884
$ws.Range("N34").Value = -8733.884
$ws.Range("H50").Value = 19999
$ws.Range("J50").Value = 19999
$ws.Range("L50").Value = 19999
$ws.Range("N50").Value = -21249
$ws.Range("H59").Value = 20165.834
$ws.Range("J59").Value = 20165.834
$ws.Range("L59").Value = 20165.834
$ws.Range("N59").Value = -22455.834
$ws.Range("H60").Value = 16173.454
$ws.Range("I60").Value = 1833.3334
$ws.Range("K60").Value = 1833.3334
$ws.Range("M60").Value = -1322.3334
$ws.Range("H68").Value = 22881.916
$ws.Range("J68").Value = 22881.916
$ws.Range("L68").Value = 22881.916
$ws.Range("N68").Value = -24379.916
$ws.Range("H71").Value = 22881.916
$ws.Range("J71").Value = 22881.916
$ws.Range("L71").Value = 68645.74800000001
$ws.Range("N71").Value = -76133.74800000001
$ws.Range("H74").Value = 19635.818
$ws.Range("J74").Value = 19635.818
$ws.Range("L74").Value = 19635.818
$ws.Range("N74").Value = -21383.818
$ws.Range("H77").Value = 19635.818
$ws.Range("J77").Value = 19635.818
$ws.Range("L77").Value = 58907.454
$ws.Range("N77").Value = -67643.454
$ws.Range("H113").Value = 2198.2104
$ws.Range("I113").Value = 2349.0833
$ws.Range("K113").Value = 2349.0833
$ws.Range("M113").Value = -179.0832999999998
$ws.Range("H137").Value = 45633.332
$ws.Range("J137").Value = 45633.332
$ws.Range("L137").Value = 45633.332
$ws.Range("N137").Value = -55833.332
$ws.Range("H141").Value = 77405.875
$ws.Range("J141").Value = 77185.5
$ws.Range("L141").Value = 77185.5
$ws.Range("N141").Value = -87545.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1954.1538
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()  # was -2829
$ws.Range("H81").Value = 5300
$ws.Range("I81").Value = 1750
$ws.Range("J81").Value = 19500
$ws.Range("K81").Value = 5250
$ws.Range("L81").Value = 58500
$ws.Range("M81").Value = -4127
$ws.Range("N81").Value = -60746
$ws.Range("H84").Value = 5300
$ws.Range("I84").Value = 1750
$ws.Range("J84").Value = 19500
$ws.Range("K84").Value = 15750
$ws.Range("L84").Value = 175500
$ws.Range("M84").Value = -10134
$ws.Range("N84").Value = -186732
$ws.Range("H87").Value = 1832.4286
$ws.Range("I87").Value = 971.1667
$ws.Range("J87").Value = 7000
$ws.Range("K87").Value = 2913.5001
$ws.Range("L87").Value = 21000
$ws.Range("M87").Value = -1665.5001
$ws.Range("N87").Value = -23496
$ws.Range("H88").Value = 5672.6665
$ws.Range("J88").Value = 5672.6665
$ws.Range("L88").Value = 17017.9995
$ws.Range("N88").Value = -17873.9995
$ws.Range("H90").Value = 1832.4286
$ws.Range("I90").Value = 971.1667
$ws.Range("J90").Value = 7000
$ws.Range("K90").Value = 8740.5003
$ws.Range("L90").Value = 63000
$ws.Range("M90").Value = -2500.5003
$ws.Range("N90").Value = -75480
$ws.Range("H91").Value = 5672.6665
$ws.Range("J91").Value = 5672.6665
$ws.Range("L91").Value = 17017.9995
$ws.Range("N91").Value = -19981.9995
$ws.Range("H107").Value = 19608760
$ws.Range("I107").Value = 187.92857
$ws.Range("J107").Value = 27028220
$ws.Range("K107").Value = 563.78571
$ws.Range("L107").Value = 81084660
$ws.Range("M107").Value = 1356.21429
$ws.Range("N107").Value = -81088500
$ws.Range("H109").Value = 2644.353
$ws.Range("I109").Value = 890.8
$ws.Range("J109").Value = 3375
$ws.Range("K109").Value = 2672.4
$ws.Range("L109").Value = 10125
$ws.Range("M109").Value = -1632.4
$ws.Range("N109").Value = -12205
$ws.Range("H112").Value = 8214.333000000001
$ws.Range("I112").Value = 975
$ws.Range("J112").Value = 8872.454
$ws.Range("K112").Value = 2925
$ws.Range("L112").Value = 26617.362
$ws.Range("M112").Value = -1817
$ws.Range("N112").Value = -28833.362
$ws.Range("H137").Value = 10454993
$ws.Range("I137").Value = 10042.643
$ws.Range("J137").Value = 18578844
$ws.Range("K137").Value = 30127.929
$ws.Range("L137").Value = 55736532
$ws.Range("M137").Value = -25027.929
$ws.Range("N137").Value = -55746732

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 80003
$ws.Range("J7").Value = 80003
$ws.Range("L7").Value = 80003
$ws.Range("N7").Value = -80227
$ws.Range("H8").Value = 80003
$ws.Range("J8").Value = 80003
$ws.Range("L8").Value = 80003
$ws.Range("N8").Value = -80281
$ws.Range("H38").Value = 10000
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10926
$ws.Range("H70").Value = 4815.1816
$ws.Range("I70").Value = 4776.2085
$ws.Range("J70").Value = 4919.1113
$ws.Range("K70").Value = 4776.2085
$ws.Range("L70").Value = 4919.1113
$ws.Range("M70").Value = -4506.2085
$ws.Range("N70").Value = -5459.1113
$ws.Range("H73").Value = 4815.1816
$ws.Range("I73").Value = 4776.2085
$ws.Range("J73").Value = 4919.1113
$ws.Range("K73").Value = 4776.2085
$ws.Range("L73").Value = 4919.1113
$ws.Range("M73").Value = -3840.2085
$ws.Range("N73").Value = -6791.1113
